# Add "Theory %" / "Practical %" rows to Metadata, and add a new "Practical"
# worksheet (a copy of the "Assg" marks sheet) positioned between "Assg" and
# "ESEM".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: append Theory % / Practical % rows
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A6").Value = "Theory %"
$meta.Range("B6").Value = 75
$meta.Range("A7").Value = "Practical %"
$meta.Range("B7").Value = 25
$meta.Rows.Item(6).RowHeight = 13.8
$meta.Rows.Item(7).RowHeight = 13.8
$meta.Range("C8").Select()

# ---------------------------------------------------------------------
# 2. Duplicate "Assg" into a new sheet named "Practical", inserted right
#    after "Assg" (i.e. right before "ESEM"), matching the marks that were
#    copied over from the Assg sheet.
# ---------------------------------------------------------------------
$assg = $wb.Worksheets.Item("Assg")
$assg.Copy($null, $assg)
$practical = $wb.Worksheets.Item($assg.Index + 1)
$practical.Name = "Practical"
$practical.Range("A6:G36").RowHeight = 13.8
$practical.Range("B7").Select()

# Keep the selection on the original Assg sheet in sync with the edit.
$assg.Range("B7").Select()

# ---------------------------------------------------------------------
# 3. The old "ESEM" sheet keeps its data untouched but now sits after the
#    new "Practical" sheet; just move the active selection.
# ---------------------------------------------------------------------
$esem = $wb.Worksheets.Item("ESEM")
$esem.Range("A6").Select()

# Leave the active sheet on Metadata, matching the workbook's original
# "tabSelected" sheet.
$meta.Activate()
